$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 366 (shifts old rows 366:371 down to 369:374,
# carrying their existing content/formatting along).
$ws.Rows("366:368").Insert()

# Row 366 - new weekly data (Packham's Triumph, Especial)
$ws.Range("A366").Value = 8
$ws.Range("B366").Value = "Terminal La Palmera de La Serena"
$ws.Range("C366").Value = "Coquimbo"
$ws.Range("D366").Value = 44509
$ws.Range("E366").Value = 4
$ws.Range("F366").Value = "Fruta"
$ws.Range("G366").Value = 100104
$ws.Range("H366").Value = "Frutos de pepita"
$ws.Range("I366").Value = 100104005
$ws.Range("J366").Value = "Pera"
$ws.Range("K366").Value = "Packham's Triumph"
$ws.Range("L366").Value = "Especial"
$ws.Range("M366").Value = 24
$ws.Range("N366").Value = 285000
$ws.Range("O366").Value = 290000
$ws.Range("P366").Value = 287500
$ws.Range("Q366").Value = "$/bins (450 kilos)"
$ws.Range("R366").Value = "Región de O'Higgins"
$ws.Range("S366").Value = 639
$ws.Range("T366").Value = 450

# Row 367 - new weekly data (Packham's Triumph, Primera)
$ws.Range("A367").Value = 8
$ws.Range("B367").Value = "Terminal La Palmera de La Serena"
$ws.Range("C367").Value = "Coquimbo"
$ws.Range("D367").Value = 44509
$ws.Range("E367").Value = 4
$ws.Range("F367").Value = "Fruta"
$ws.Range("G367").Value = 100104
$ws.Range("H367").Value = "Frutos de pepita"
$ws.Range("I367").Value = 100104005
$ws.Range("J367").Value = "Pera"
$ws.Range("K367").Value = "Packham's Triumph"
$ws.Range("L367").Value = "Primera"
$ws.Range("M367").Value = 20
$ws.Range("N367").Value = 255000
$ws.Range("O367").Value = 260000
$ws.Range("P367").Value = 257500
$ws.Range("Q367").Value = "$/bins (450 kilos)"
$ws.Range("R367").Value = "Región de O'Higgins"
$ws.Range("S367").Value = 572
$ws.Range("T367").Value = 450

# Row 368 - new weekly data (Packham's Triumph, Segunda)
$ws.Range("A368").Value = 8
$ws.Range("B368").Value = "Terminal La Palmera de La Serena"
$ws.Range("C368").Value = "Coquimbo"
$ws.Range("D368").Value = 44509
$ws.Range("E368").Value = 4
$ws.Range("F368").Value = "Fruta"
$ws.Range("G368").Value = 100104
$ws.Range("H368").Value = "Frutos de pepita"
$ws.Range("I368").Value = 100104005
$ws.Range("J368").Value = "Pera"
$ws.Range("K368").Value = "Packham's Triumph"
$ws.Range("L368").Value = "Segunda"
$ws.Range("M368").Value = 18
$ws.Range("N368").Value = 235000
$ws.Range("O368").Value = 240000
$ws.Range("P368").Value = 237500
$ws.Range("Q368").Value = "$/bins (450 kilos)"
$ws.Range("R368").Value = "Región de O'Higgins"
$ws.Range("S368").Value = 528
$ws.Range("T368").Value = 450
